$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new bibliography entry in row 17 (in-text citation + full reference)
$ws.Range("A17").Value = "(Oliver & Snowden, 2015)"
$ws.Range("B17").Value = "Oliver, J., & Snowden, E. [Last Week Tonight]. (2015). Last Week Tonight with John Oliver: Edward Snowden on Passwords. Retrieved May 6, 2015, from https://www.youtube.com/watch?v=yzGzB-yYKcc"

# Match the left-aligned style used by the rest of column B
$ws.Range("B17").HorizontalAlignment = -4131

# Update the active cell selection like in the saved workbook
$ws.Range("B20").Select()
